# Update countries & provincias Spain
# Applies the daily COVID data refresh captured in the diff:
#  - "Nepal" overtakes Portugal/Etiopia/Costa Rica/Venezuela in the ranking
#    (rows 51-55 shift down, row 51 gets Nepal's new totals)
#  - "Timor Oriental" and "Santa Lucia" swap places (rows 206/207, tied data)
#  - A batch of per-country case/death counts are refreshed
#  - The "last updated" timestamp moves from 12:23 to 13:40

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 13:40"

# --- straightforward numeric refresh: country stays put, only B..H move ---
$ws.Range("B16").Value = 446448
$ws.Range("C16").Value = 3362
$ws.Range("D16").Value = 374170
$ws.Range("E16").Value = 46689
$ws.Range("G16").Value = 195
$ws.Range("H16").Value = 25589

$ws.Range("B18").Value = 359148
$ws.Range("C18").Value = 1275
$ws.Range("D18").Value = 270491
$ws.Range("E18").Value = 83496
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 5161

$ws.Range("B42").Value = 97450
$ws.Range("C42").Value = 1543
$ws.Range("D42").Value = 87801
$ws.Range("E42").Value = 8740
$ws.Range("G42").Value = 24
$ws.Range("H42").Value = 909

$ws.Range("B49").Value = 77609
$ws.Range("C49").Value = 320
$ws.Range("D49").Value = 74120
$ws.Range("E49").Value = 2671
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 818

$ws.Range("B68").Value = 40023
$ws.Range("C68").Value = 128
$ws.Range("D68").Value = 37655
$ws.Range("E68").Value = 1782
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 586

$ws.Range("B88").Value = 16285
$ws.Range("C88").Value = 28
$ws.Range("E88").Value = 1134

$ws.Range("B91").Value = 14909
$ws.Range("C91").Value = 40
$ws.Range("D91").Value = 12113
$ws.Range("E91").Value = 2488
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 308

$ws.Range("B139").Value = 3564
$ws.Range("C139").Value = 9
$ws.Range("D139").Value = 2061
$ws.Range("E139").Value = 1393

$ws.Range("D142").Value = 3208
$ws.Range("E142").Value = 128

$ws.Range("B145").Value = 2979
$ws.Range("C145").Value = 21
$ws.Range("D145").Value = 2358
$ws.Range("E145").Value = 590

$ws.Range("B168").Value = 1074
$ws.Range("C168").Value = 5
$ws.Range("E168").Value = 40

# --- Nepal climbs above Portugal / Etiopia / Costa Rica / Venezuela ---
# Row 51 now carries Nepal's own (updated) totals; rows 52-55 inherit the
# values the row above used to hold (Portugal, Etiopia, Costa Rica,
# Venezuela respectively), each shifted down one rank; row 56 (Barein)
# is untouched.
$ws.Range("A51").Value = "Nepal"
$ws.Range("B51").Value = 73394
$ws.Range("C51").Value = 1573
$ws.Range("D51").Value = 53898
$ws.Range("E51").Value = 19019
$ws.Range("G51").Value = 10
$ws.Range("H51").Value = 477

$ws.Range("A52").Value = "Portugal"
$ws.Range("B52").Value = 72939
$ws.Range("D52").Value = 47380
$ws.Range("E52").Value = 23615
$ws.Range("H52").Value = 1944

$ws.Range("A53").Value = "Etiopia"
$ws.Range("B53").Value = 72700
$ws.Range("D53").Value = 30029
$ws.Range("E53").Value = 41506
$ws.Range("H53").Value = 1165

$ws.Range("A54").Value = "Costa Rica"
$ws.Range("B54").Value = 72049
$ws.Range("D54").Value = 27760
$ws.Range("E54").Value = 43461
$ws.Range("H54").Value = 828

$ws.Range("A55").Value = "Venezuela"
$ws.Range("B55").Value = 71940
$ws.Range("D55").Value = 61528
$ws.Range("E55").Value = 9812
$ws.Range("H55").Value = 600

# --- Timor Oriental / Santa Lucia swap ranks (identical totals, so only
#     the country labels need to trade places) ---
$ws.Range("A206").Value = "Timor Oriental"
$ws.Range("A207").Value = "Santa Lucia"
